$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Change Parameters" (row 5) / "Change Power" (row 6) flags flipped from Yes to No
$ws.Range("B5").Value = "No"
$ws.Range("B6").Value = "No"

# Selection moved from B5 to B7
$ws.Range("B7").Select()
